$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.529.38'
$ws.Range('E2').Value = '  +2.35%  '
$ws.Range('D3').Value = '1.857.56'
$ws.Range('E3').Value = '  +1.51%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.57'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6950'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.12%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3073'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07680'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.16%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.61'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.15%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07768'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.60%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.869.32'
$ws.Range('E12').Value = '  +2.21%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.156'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.62%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6947'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.21%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '90.98'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.94%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.307'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.98%  '
$ws.Range('D17').Value = '29.526.73'
$ws.Range('E17').Value = '  +2.35%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008294'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.16%  '
$ws.Range('D19').Value = '2.103.06'
$ws.Range('E19').Value = '  +1.51%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '237.33'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.67%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.75'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.64%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.000'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.630'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.51%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.001'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1488'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.13%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.907'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.56%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '160.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.26%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.26'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.30%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.536'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.44%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.248'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.87%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.147'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.27%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.215'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.39%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05301'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7793'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.47%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.880'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.54%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.148'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.11%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.683'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.42%  '
$ws.Range('D38').Value = '1.318.79'
$ws.Range('E38').Value = '  +8.18%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01872'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.54%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.723'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.93%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9459'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.75%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '105.99'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.93%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.757'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.02%  '
$ws.Range('E44').Value = '  +0.10%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '9.760'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.58%  '
$ws.Range('E46').Value = '  +1.93%  '
$ws.Range('D47').Value = '2.006.52'
$ws.Range('E47').Value = '  +1.63%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.786'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.43%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '62.83'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.35%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05961'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.71%  '
